# feat: add 2022-Q3 data
#
# Inserts a new "2022-Q3" sheet (cloned from "2022-Q2") populated with the
# Q3 fund-holding data, updates the "总计" (totals) summary sheet with a new
# leading row for 2022-Q3, and shifts the existing 2022-Q2 / 2022-Q1 rows
# down to make room.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "总计" sheet: insert the 2022-Q3 total as the new first data row,
#    pushing 2022-Q2 and 2022-Q1 down by one row each. Formatting for the
#    "#" column (A) is carried over via Copy so the shared cell style
#    (centered/bold/bordered) is reused rather than recreated.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

# Row 3 (2022-Q1) -> Row 4
$totals.Range("A3").Copy($totals.Range("A4"))
$totals.Range("B3:D3").Copy($totals.Range("B4:D4"))
$totals.Range("A4").Value = 2
$totals.Range("B4").Value = "2022-Q1"
$totals.Range("C4").Value = 2
$totals.Range("D4").Value = 0.01

# Row 2 (2022-Q2) -> Row 3
$totals.Range("A2").Copy($totals.Range("A3"))
$totals.Range("B2:D2").Copy($totals.Range("B3:D3"))
$totals.Range("A3").Value = 1
$totals.Range("B3").Value = "2022-Q2"
$totals.Range("C3").Value = 5
$totals.Range("D3").Value = 0.13

# Row 2 becomes the new 2022-Q3 total
$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q3"
$totals.Range("C2").Value = 1
$totals.Range("D2").Value = 0.03

# ---------------------------------------------------------------------
# 2. Build the "2022-Q3" detail sheet. Cloning "2022-Q2" keeps every style
#    (header shading/border, "#" column formatting, column widths, etc.)
#    identical, then we trim it down to the single Q3 holding and overwrite
#    the values.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($null, $totals)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# Drop the old Q2 sample rows below the header + first data row, leaving
# just row 2 to be overwritten with the Q3 holding.
$q3.Range("A3:H6").Delete(-4162)

$q3.Range("A2").Value = 0

# Text-ish numeric strings (fund code / percentages) need to stay text, as
# in the sibling quarter sheets - flip to Text format just long enough to
# assign the literal string, then clear the format override so the cell
# keeps the sheet's default (unstyled) look, matching B2/D2:G2 elsewhere.
$q3.Range("B2").NumberFormat = "@"
$q3.Range("B2").Value = "001294"
$q3.Range("B2").ClearFormats()

$q3.Range("C2").Value = "新华战略新兴产业灵活配置混合"

$q3.Range("D2").NumberFormat = "@"
$q3.Range("D2").Value = "0.99"
$q3.Range("D2").ClearFormats()

$q3.Range("E2").NumberFormat = "@"
$q3.Range("E2").Value = "93.49"
$q3.Range("E2").ClearFormats()

$q3.Range("F2").NumberFormat = "@"
$q3.Range("F2").Value = "3.24"
$q3.Range("F2").ClearFormats()

$q3.Range("G2").NumberFormat = "@"
$q3.Range("G2").Value = "0.0321"
$q3.Range("G2").ClearFormats()

$q3.Range("H2").Value = 9

# Restore "总计" as the active sheet (it was the active/selected sheet
# before this edit and the edit doesn't change that).
$totals.Activate()
